# Add new Testcafe signup API user credentials to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$users = @(
    "user_34dca657-17a8-11ec-8bdd-38fc98d48cf7@mail.com",
    "user_48012c6f-17a8-11ec-9bb0-38fc98d48cf7@mail.com",
    "user_4b864083-17a8-11ec-a0e7-38fc98d48cf7@mail.com",
    "user_4daf9901-17a8-11ec-bad2-38fc98d48cf7@mail.com",
    "user_0415c716-17a9-11ec-8f1b-38fc98d48cf7@mail.com",
    "user_820d387e-17a9-11ec-82f4-38fc98d48cf7@mail.com",
    "user_d97adb5b-17a9-11ec-ba2e-38fc98d48cf7@mail.com",
    "user_26a78e36-17aa-11ec-bd2d-38fc98d48cf7@mail.com",
    "user_56d5978c-17ab-11ec-aacc-38fc98d48cf7@mail.com",
    "user_7a315277-17ae-11ec-ad3a-38fc98d48cf7@mail.com"
)

$password = "Asdfgh123!"

$row = 3
foreach ($user in $users) {
    $ws.Cells.Item($row, 1).Value = $user
    $ws.Cells.Item($row, 2).Value = $password
    $row = $row + 1
}
